# "Processed Results - nexus 5x until gps"
#
# The CPU benchmark run for the Nexus 5X was re-measured, so the raw
# per-run energy samples in column B of "Blad1" change for all three
# frequency blocks:
#   - Low frequency samples:    B2:B31   (mirrored again at B98:B127)
#   - Medium frequency samples: B34:B63  (mirrored again at B128:B157)
#   - High frequency samples:   B66:B95  (mirrored again at B158:B187)
# Column A (run/frequency labels) and all summary formulas in columns D/E
# (AVERAGE/MIN/MAX/QUARTILE/STDEV) are left untouched - they recompute
# automatically from the new B values, as do the charts built on top of
# them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lowVals  = @(122.01944399999999, 129.537972, 123.585804, 124.838892, 124.055712, 126.71852399999899, 127.65834, 123.115895999999, 124.21234800000001, 124.055712, 122.95926, 119.82653999999999, 124.055712, 123.89907599999999, 124.682256, 128.75479200000001, 123.115895999999, 123.429167999999, 124.368984, 124.52562, 123.115895999999, 123.429167999999, 124.055712, 127.188431999999, 122.489352, 123.89907599999999, 122.489352, 123.585804, 119.043359999999, 123.429167999999)
$medVals  = @(127.814976, 134.70695999999899, 128.59815599999999, 129.537972, 127.97161199999999, 129.537972, 128.12824800000001, 128.59815599999999, 128.44152, 128.59815599999999, 129.537972, 124.838892, 136.89986399999901, 128.59815599999999, 129.22469999999899, 128.12824800000001, 128.28488399999901, 127.814976, 129.22469999999899, 128.75479200000001, 132.98396399999999, 128.44152, 129.537972, 131.104332, 130.32115200000001, 134.23705200000001, 129.06806399999999, 133.297236, 124.52562, 129.22469999999899)
$highVals = @(144.26175599999999, 145.82811599999999, 144.88829999999999, 145.35820799999999, 147.394476, 144.26175599999999, 144.105119999999, 144.26175599999999, 143.635212, 146.92456799999999, 145.67148, 142.69539599999999, 145.201572, 144.575028, 145.35820799999999, 144.73166399999999, 146.14138800000001, 150.213923999999, 144.26175599999999, 143.635212, 144.26175599999999, 143.79184799999999, 143.635212, 145.201572, 145.04493600000001, 145.04493600000001, 144.105119999999, 143.478576, 144.26175599999999, 144.41839200000001)

for ($i = 0; $i -lt $lowVals.Length; $i++) {
    $ws.Cells.Item(2 + $i, 2).Value = $lowVals[$i]
    $ws.Cells.Item(98 + $i, 2).Value = $lowVals[$i]
}

for ($i = 0; $i -lt $medVals.Length; $i++) {
    $ws.Cells.Item(34 + $i, 2).Value = $medVals[$i]
    $ws.Cells.Item(128 + $i, 2).Value = $medVals[$i]
}

for ($i = 0; $i -lt $highVals.Length; $i++) {
    $ws.Cells.Item(66 + $i, 2).Value = $highVals[$i]
    $ws.Cells.Item(158 + $i, 2).Value = $highVals[$i]
}

# Make sure every dependent formula (averages, min/max, quartiles, stdev)
# and chart picks up the new source values.
$excel.CalculateFullRebuild()

# Leave the sheet scrolled/selected where the author ended up (row 85 area,
# cell I95) instead of the original E48.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 85
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I95").Select()
